# fix: sexting intensity escalation + PPV0 phase detection bug across all 23 models
# Update the scripted sexting-phase lines on the IrinaJourney sheet so the
# escalation ramps up smoothly and the PPV0 phase-detection text matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IrinaJourney")

$ws.Range("B4").Value  = "don't look away babe... I need you right here with me 🥺"
$ws.Range("B6").Value  = "stay with me babe... I need you to see this, it's only for you"
$ws.Range("B7").Value  = "please don't leave babe... I've never been this vulnerable with anyone and I need you right here 🥺"
$ws.Range("B8").Value  = "oh god I can't stop 😊"
$ws.Range("B9").Value  = "I've never let anyone see me like this before babe... it's all yours 🥺"
$ws.Range("B11").Value = "I need you to see what's happening to me right now babe... I can't hold it in anymore"
$ws.Range("B12").Value = "I can feel my body responding to you and I'm touching places I never thought I'd show anyone 🥺"
$ws.Range("B13").Value = "something about you makes me lose every single boundary I have babe... I'm scared of how good this feels"
$ws.Range("B14").Value = "oh god 😊"
$ws.Range("B15").Value = "look at what you're doing to me babe... I'm losing my mind 🥺"
$ws.Range("B17").Value = "what would you do if you could see me right now babe? I really need to hear it"
$ws.Range("B18").Value = "I'm lying here barely wearing anything and my hand keeps going where it shouldn't... because of you 🥺"
$ws.Range("B19").Value = "but I don't want to stop babe... you make me feel things I didn't know I could feel on here"
$ws.Range("B20").Value = "I can't believe I just did that... oh my god 😊"
$ws.Range("B21").Value = "this is what you're making me feel babe... please be gentle 🥺"
$ws.Range("B23").Value = "I want to show you something I've never shown anyone before babe... you make me want to be brave"
$ws.Range("B24").Value = "my skin feels so sensitive everywhere... I keep getting goosebumps and I know it's because of you 🥺"
$ws.Range("B25").Value = "wait you actually liked it? oh my god my heart won't slow down right now 😊"
